$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "vendredi" -> "mardi" (shared string used in B2, B5, B8)
$ws.Range("B2").Value = "mardi"
$ws.Range("B5").Value = "mardi"
$ws.Range("B8").Value = "mardi"

# Update the dates (same day-of-month, shifted to 2026, landing on Tuesday)
$ws.Range("A2").Value = 46063
$ws.Range("A5").Value = 46077
$ws.Range("A8").Value = 46091
